$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting of scraped floating point numbers (comma decimal -> dot decimal,
# drop thousands-separator dots) and normalize stray commas/periods in a few
# "Razon social" entries that were mangled by the same scrape bug.

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '660.00'
$ws.Range("H2").Style = "Normal"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '1208.00'
$ws.Range("H3").Style = "Normal"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = '34.44'
$ws.Range("H4").Style = "Normal"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '570.00'
$ws.Range("H5").Style = "Normal"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '1810.00'
$ws.Range("H6").Style = "Normal"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '34.10'
$ws.Range("H7").Style = "Normal"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '2476.56'
$ws.Range("H8").Style = "Normal"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '770.00'
$ws.Range("H9").Style = "Normal"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '1326.76'
$ws.Range("H10").Style = "Normal"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '70.00'
$ws.Range("H11").Style = "Normal"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '12760.00'
$ws.Range("H12").Style = "Normal"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '25900.00'
$ws.Range("H13").Style = "Normal"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '83321.75'
$ws.Range("H14").Style = "Normal"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '1290.50'
$ws.Range("H15").Style = "Normal"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '22213.63'
$ws.Range("H16").Style = "Normal"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = '1509.46'
$ws.Range("H17").Style = "Normal"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = '6031.66'
$ws.Range("H18").Style = "Normal"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = '6447.46'
$ws.Range("H19").Style = "Normal"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '650.00'
$ws.Range("H20").Style = "Normal"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '1260.00'
$ws.Range("H21").Style = "Normal"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = '4005.00'
$ws.Range("H22").Style = "Normal"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '200.35'
$ws.Range("H23").Style = "Normal"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = '14.00'
$ws.Range("H24").Style = "Normal"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = '792.00'
$ws.Range("H25").Style = "Normal"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '9153.00'
$ws.Range("H26").Style = "Normal"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = '258.85'
$ws.Range("H27").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E28").Style = "Normal"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F28").Style = "Normal"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '111.10'
$ws.Range("H28").Style = "Normal"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = '93.42'
$ws.Range("H29").Style = "Normal"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '6775.38'
$ws.Range("H30").Style = "Normal"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '210.00'
$ws.Range("H31").Style = "Normal"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '2680.00'
$ws.Range("H32").Style = "Normal"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '2419.48'
$ws.Range("H33").Style = "Normal"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = '20.00'
$ws.Range("H34").Style = "Normal"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '1784.06'
$ws.Range("H35").Style = "Normal"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '319.51'
$ws.Range("H36").Style = "Normal"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = '29437.36'
$ws.Range("H37").Style = "Normal"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = '12527.30'
$ws.Range("H38").Style = "Normal"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = '9144.00'
$ws.Range("H39").Style = "Normal"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '300.00'
$ws.Range("H40").Style = "Normal"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '606.20'
$ws.Range("H41").Style = "Normal"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = '4327.29'
$ws.Range("H42").Style = "Normal"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = '228.80'
$ws.Range("H43").Style = "Normal"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = '9358.00'
$ws.Range("H44").Style = "Normal"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '340.92'
$ws.Range("H45").Style = "Normal"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '678.00'
$ws.Range("H46").Style = "Normal"
$ws.Range("H47").NumberFormat = "@"
$ws.Range("H47").Value = '971.00'
$ws.Range("H47").Style = "Normal"
$ws.Range("H48").NumberFormat = "@"
$ws.Range("H48").Value = '2416.46'
$ws.Range("H48").Style = "Normal"
$ws.Range("H49").NumberFormat = "@"
$ws.Range("H49").Value = '1440.50'
$ws.Range("H49").Style = "Normal"
$ws.Range("H50").NumberFormat = "@"
$ws.Range("H50").Value = '11328.00'
$ws.Range("H50").Style = "Normal"
$ws.Range("H51").NumberFormat = "@"
$ws.Range("H51").Value = '1300.00'
$ws.Range("H51").Style = "Normal"
$ws.Range("H52").NumberFormat = "@"
$ws.Range("H52").Value = '100.00'
$ws.Range("H52").Style = "Normal"
$ws.Range("H53").NumberFormat = "@"
$ws.Range("H53").Value = '4674.67'
$ws.Range("H53").Style = "Normal"
$ws.Range("E54").NumberFormat = "@"
$ws.Range("E54").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("E54").Style = "Normal"
$ws.Range("F54").NumberFormat = "@"
$ws.Range("F54").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("F54").Style = "Normal"
$ws.Range("H54").NumberFormat = "@"
$ws.Range("H54").Value = '200.00'
$ws.Range("H54").Style = "Normal"
$ws.Range("H55").NumberFormat = "@"
$ws.Range("H55").Value = '1509.10'
$ws.Range("H55").Style = "Normal"
$ws.Range("H56").NumberFormat = "@"
$ws.Range("H56").Value = '10595.20'
$ws.Range("H56").Style = "Normal"
$ws.Range("H57").NumberFormat = "@"
$ws.Range("H57").Value = '12250.00'
$ws.Range("H57").Style = "Normal"
$ws.Range("H58").NumberFormat = "@"
$ws.Range("H58").Value = '8505.00'
$ws.Range("H58").Style = "Normal"
$ws.Range("H59").NumberFormat = "@"
$ws.Range("H59").Value = '4176.00'
$ws.Range("H59").Style = "Normal"
$ws.Range("H60").NumberFormat = "@"
$ws.Range("H60").Value = '1301.80'
$ws.Range("H60").Style = "Normal"
$ws.Range("E61").NumberFormat = "@"
$ws.Range("E61").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E61").Style = "Normal"
$ws.Range("H61").NumberFormat = "@"
$ws.Range("H61").Value = '4038.00'
$ws.Range("H61").Style = "Normal"
$ws.Range("E62").NumberFormat = "@"
$ws.Range("E62").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E62").Style = "Normal"
$ws.Range("F62").NumberFormat = "@"
$ws.Range("F62").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F62").Style = "Normal"
$ws.Range("H62").NumberFormat = "@"
$ws.Range("H62").Value = '3479.54'
$ws.Range("H62").Style = "Normal"
$ws.Range("H63").NumberFormat = "@"
$ws.Range("H63").Value = '370.00'
$ws.Range("H63").Style = "Normal"
$ws.Range("E64").NumberFormat = "@"
$ws.Range("E64").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E64").Style = "Normal"
$ws.Range("H64").NumberFormat = "@"
$ws.Range("H64").Value = '208.00'
$ws.Range("H64").Style = "Normal"
$ws.Range("H65").NumberFormat = "@"
$ws.Range("H65").Value = '5814.00'
$ws.Range("H65").Style = "Normal"
$ws.Range("H66").NumberFormat = "@"
$ws.Range("H66").Value = '95670.00'
$ws.Range("H66").Style = "Normal"
$ws.Range("H67").NumberFormat = "@"
$ws.Range("H67").Value = '330.40'
$ws.Range("H67").Style = "Normal"
$ws.Range("H68").NumberFormat = "@"
$ws.Range("H68").Value = '325.48'
$ws.Range("H68").Style = "Normal"
$ws.Range("H69").NumberFormat = "@"
$ws.Range("H69").Value = '623.00'
$ws.Range("H69").Style = "Normal"
$ws.Range("H70").NumberFormat = "@"
$ws.Range("H70").Value = '408.00'
$ws.Range("H70").Style = "Normal"
$ws.Range("H71").NumberFormat = "@"
$ws.Range("H71").Value = '961.30'
$ws.Range("H71").Style = "Normal"
$ws.Range("H72").NumberFormat = "@"
$ws.Range("H72").Value = '12.86'
$ws.Range("H72").Style = "Normal"
$ws.Range("H73").NumberFormat = "@"
$ws.Range("H73").Value = '17105.55'
$ws.Range("H73").Style = "Normal"
$ws.Range("H74").NumberFormat = "@"
$ws.Range("H74").Value = '84.00'
$ws.Range("H74").Style = "Normal"
$ws.Range("H75").NumberFormat = "@"
$ws.Range("H75").Value = '16.00'
$ws.Range("H75").Style = "Normal"
$ws.Range("E76").NumberFormat = "@"
$ws.Range("E76").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E76").Style = "Normal"
$ws.Range("H76").NumberFormat = "@"
$ws.Range("H76").Value = '8846.90'
$ws.Range("H76").Style = "Normal"
$ws.Range("H77").NumberFormat = "@"
$ws.Range("H77").Value = '20969.00'
$ws.Range("H77").Style = "Normal"
$ws.Range("E78").NumberFormat = "@"
$ws.Range("E78").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E78").Style = "Normal"
$ws.Range("F78").NumberFormat = "@"
$ws.Range("F78").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F78").Style = "Normal"
$ws.Range("H78").NumberFormat = "@"
$ws.Range("H78").Value = '4019.20'
$ws.Range("H78").Style = "Normal"
$ws.Range("H79").NumberFormat = "@"
$ws.Range("H79").Value = '16696.80'
$ws.Range("H79").Style = "Normal"
$ws.Range("H80").NumberFormat = "@"
$ws.Range("H80").Value = '4055.39'
$ws.Range("H80").Style = "Normal"
$ws.Range("H81").NumberFormat = "@"
$ws.Range("H81").Value = '17.16'
$ws.Range("H81").Style = "Normal"
$ws.Range("H82").NumberFormat = "@"
$ws.Range("H82").Value = '1183.70'
$ws.Range("H82").Style = "Normal"
$ws.Range("H83").NumberFormat = "@"
$ws.Range("H83").Value = '1617.20'
$ws.Range("H83").Style = "Normal"
$ws.Range("H84").NumberFormat = "@"
$ws.Range("H84").Value = '52.49'
$ws.Range("H84").Style = "Normal"
$ws.Range("H85").NumberFormat = "@"
$ws.Range("H85").Value = '118.25'
$ws.Range("H85").Style = "Normal"
$ws.Range("H86").NumberFormat = "@"
$ws.Range("H86").Value = '37.60'
$ws.Range("H86").Style = "Normal"
$ws.Range("H87").NumberFormat = "@"
$ws.Range("H87").Value = '35.00'
$ws.Range("H87").Style = "Normal"
$ws.Range("H88").NumberFormat = "@"
$ws.Range("H88").Value = '26920.00'
$ws.Range("H88").Style = "Normal"
$ws.Range("H89").NumberFormat = "@"
$ws.Range("H89").Value = '860.00'
$ws.Range("H89").Style = "Normal"
$ws.Range("H90").NumberFormat = "@"
$ws.Range("H90").Value = '149.16'
$ws.Range("H90").Style = "Normal"
$ws.Range("H91").NumberFormat = "@"
$ws.Range("H91").Value = '185.07'
$ws.Range("H91").Style = "Normal"
$ws.Range("H92").NumberFormat = "@"
$ws.Range("H92").Value = '713.50'
$ws.Range("H92").Style = "Normal"
$ws.Range("H93").NumberFormat = "@"
$ws.Range("H93").Value = '787.93'
$ws.Range("H93").Style = "Normal"
$ws.Range("H94").NumberFormat = "@"
$ws.Range("H94").Value = '250.00'
$ws.Range("H94").Style = "Normal"
$ws.Range("H95").NumberFormat = "@"
$ws.Range("H95").Value = '500.00'
$ws.Range("H95").Style = "Normal"
$ws.Range("H96").NumberFormat = "@"
$ws.Range("H96").Value = '3500.00'
$ws.Range("H96").Style = "Normal"
$ws.Range("H97").NumberFormat = "@"
$ws.Range("H97").Value = '2504.70'
$ws.Range("H97").Style = "Normal"
$ws.Range("H98").NumberFormat = "@"
$ws.Range("H98").Value = '290.00'
$ws.Range("H98").Style = "Normal"
$ws.Range("H99").NumberFormat = "@"
$ws.Range("H99").Value = '1815.00'
$ws.Range("H99").Style = "Normal"
$ws.Range("H100").NumberFormat = "@"
$ws.Range("H100").Value = '250.00'
$ws.Range("H100").Style = "Normal"
$ws.Range("H101").NumberFormat = "@"
$ws.Range("H101").Value = '760.00'
$ws.Range("H101").Style = "Normal"
$ws.Range("H102").NumberFormat = "@"
$ws.Range("H102").Value = '7127.25'
$ws.Range("H102").Style = "Normal"
$ws.Range("H103").NumberFormat = "@"
$ws.Range("H103").Value = '600.00'
$ws.Range("H103").Style = "Normal"
$ws.Range("H104").NumberFormat = "@"
$ws.Range("H104").Value = '350.00'
$ws.Range("H104").Style = "Normal"
$ws.Range("H105").NumberFormat = "@"
$ws.Range("H105").Value = '870.00'
$ws.Range("H105").Style = "Normal"
$ws.Range("H106").NumberFormat = "@"
$ws.Range("H106").Value = '120.00'
$ws.Range("H106").Style = "Normal"
$ws.Range("H107").NumberFormat = "@"
$ws.Range("H107").Value = '18075.00'
$ws.Range("H107").Style = "Normal"
$ws.Range("H108").NumberFormat = "@"
$ws.Range("H108").Value = '198.00'
$ws.Range("H108").Style = "Normal"
$ws.Range("E109").NumberFormat = "@"
$ws.Range("E109").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E109").Style = "Normal"
$ws.Range("H109").NumberFormat = "@"
$ws.Range("H109").Value = '170.00'
$ws.Range("H109").Style = "Normal"
$ws.Range("H110").NumberFormat = "@"
$ws.Range("H110").Value = '350.00'
$ws.Range("H110").Style = "Normal"
$ws.Range("H111").NumberFormat = "@"
$ws.Range("H111").Value = '1210.00'
$ws.Range("H111").Style = "Normal"
$ws.Range("H112").NumberFormat = "@"
$ws.Range("H112").Value = '1144.00'
$ws.Range("H112").Style = "Normal"
$ws.Range("H113").NumberFormat = "@"
$ws.Range("H113").Value = '3610.96'
$ws.Range("H113").Style = "Normal"
$ws.Range("H114").NumberFormat = "@"
$ws.Range("H114").Value = '1226.00'
$ws.Range("H114").Style = "Normal"
$ws.Range("H115").NumberFormat = "@"
$ws.Range("H115").Value = '1401.00'
$ws.Range("H115").Style = "Normal"
$ws.Range("H116").NumberFormat = "@"
$ws.Range("H116").Value = '1968.56'
$ws.Range("H116").Style = "Normal"
$ws.Range("H117").NumberFormat = "@"
$ws.Range("H117").Value = '168.00'
$ws.Range("H117").Style = "Normal"
$ws.Range("H118").NumberFormat = "@"
$ws.Range("H118").Value = '1003.10'
$ws.Range("H118").Style = "Normal"
$ws.Range("H119").NumberFormat = "@"
$ws.Range("H119").Value = '59.40'
$ws.Range("H119").Style = "Normal"
$ws.Range("H120").NumberFormat = "@"
$ws.Range("H120").Value = '404.39'
$ws.Range("H120").Style = "Normal"
$ws.Range("H121").NumberFormat = "@"
$ws.Range("H121").Value = '1809.10'
$ws.Range("H121").Style = "Normal"
$ws.Range("H122").NumberFormat = "@"
$ws.Range("H122").Value = '6536.00'
$ws.Range("H122").Style = "Normal"
$ws.Range("H123").NumberFormat = "@"
$ws.Range("H123").Value = '1758.00'
$ws.Range("H123").Style = "Normal"
$ws.Range("H124").NumberFormat = "@"
$ws.Range("H124").Value = '8174.48'
$ws.Range("H124").Style = "Normal"
$ws.Range("H125").NumberFormat = "@"
$ws.Range("H125").Value = '2685.22'
$ws.Range("H125").Style = "Normal"
$ws.Range("H126").NumberFormat = "@"
$ws.Range("H126").Value = '11200.00'
$ws.Range("H126").Style = "Normal"
$ws.Range("H127").NumberFormat = "@"
$ws.Range("H127").Value = '1127.36'
$ws.Range("H127").Style = "Normal"
$ws.Range("H128").NumberFormat = "@"
$ws.Range("H128").Value = '344598.13'
$ws.Range("H128").Style = "Normal"
$ws.Range("H129").NumberFormat = "@"
$ws.Range("H129").Value = '12800.00'
$ws.Range("H129").Style = "Normal"
$ws.Range("H130").NumberFormat = "@"
$ws.Range("H130").Value = '55500.00'
$ws.Range("H130").Style = "Normal"
$ws.Range("H131").NumberFormat = "@"
$ws.Range("H131").Value = '231255.63'
$ws.Range("H131").Style = "Normal"
$ws.Range("H132").NumberFormat = "@"
$ws.Range("H132").Value = '456.00'
$ws.Range("H132").Style = "Normal"
